# Apply "updated styles of main courses" edit.
$d = $word.ActiveDocument

# --- 1. Update the Heading4 paragraph style: drop the explicit orange color ---
$h4 = $d.Styles.Item("Heading4")
$h4.Font.Color = -16777216   # wdColorAutomatic (reset to "no explicit color")

# --- 2. Update the Heading4Char (linked character) style ---
$h4c = $d.Styles.Item("Heading4Char")
$h4c.Font.Italic = $false
$h4c.Font.Bold = $true
$h4c.Font.Color = 20736      # RGB(0, 176, 80) = 0x00B050

# --- 3. Remove direct run-level green color formatting on the three
#        Heading4 paragraphs (Ingrédients / Préparation / Viande), which
#        otherwise overrides the (now-updated) style color. ---
foreach ($p in $d.Paragraphs) {
    if ($p.Style.NameLocal -eq "Heading 4") {
        $p.Range.Font.Color = -16777216  # wdColorAutomatic
    }
}

# --- 4. Merge the "steack" run-split (proofErr spellcheck bookmarks) back
#        into a single plain run of text. ---
$d.Content.Find.Execute( `
    "Peut être servi en accompagnement de toutes les viandes: steack haché, bœuf à griller, rôtis…", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Peut être servi en accompagnement de toutes les viandes: steack haché, bœuf à griller, rôtis…", `
    2) | Out-Null
